$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-6: rotate the site names, keep B/C/D (200, Checked, Success) as-is
$ws.Range("A2").Value = "https://www.google.com"
$ws.Range("A3").Value = "https://www.sebi.gov.in"
$ws.Range("A4").Value = "https://www.surveymonkey.com"
$ws.Range("A5").Value = "https://www.owasp.org"
$ws.Range("A6").Value = "https://www.axisbank.com"

# Row 7: new unresolvable domain, not checked, with error reason; no response code
$ws.Range("A7").Value = "https://www.sahilendworldfibvweuidbuk.org"
$ws.Range("B7").ClearContents()
$ws.Range("C7").Value = "Not Checked"
$ws.Range("D7").Value = 'HTTPSConnectionPool(host=''www.sahilendworldfibvweuidbuk.org'', port=443): Max retries exceeded with url: / (Caused by NameResolutionError("<urllib3.connection.HTTPSConnection object at 0x000001CC7FEB4550>: Failed to resolve ''www.sahilendworldfibvweuidbuk.org'' ([Errno 11001] getaddrinfo failed)"))'

# Row 8 (new): rbi.org.in moves to the end of the list, fully checked
$ws.Range("A8").Value = "https://www.rbi.org.in"
$ws.Range("B8").Value = 200
$ws.Range("C8").Value = "Checked"
$ws.Range("D8").Value = "Success"
